# Insert a new data row at row 346 (pushing the existing rows 346-446 down
# to 347-447) and populate it with the new record, matching the author's
# weekly price update for "Hortaliza, Feria Lagunitas de Puerto Montt - Sandia".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 346:446 down by one row, creating a blank row 346.
$ws.Rows("346:346").Insert()

# Fill the newly inserted row with the new record's values.
$ws.Range("A346").Value = 4
$ws.Range("B346").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C346").Value = "Los Lagos"
$ws.Range("D346").Value = 45215
$ws.Range("E346").Value = 10
$ws.Range("F346").Value = 100112028
$ws.Range("G346").Value = "Sandia"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Segunda"
$ws.Range("J346").Value = 250
$ws.Range("K346").Value = 1000
$ws.Range("L346").Value = 1000
$ws.Range("M346").Value = 1000
$ws.Range("N346").Value = '$/kilo (volumen en unidades)'
$ws.Range("O346").Value = "Perú"
$ws.Range("P346").Value = 1000
$ws.Range("Q346").Value = 1
$ws.Range("R346").Value = "Hortaliza"
